# CA Grades.docx - "Added end of year stuff"
#
# Changes applied:
#  1. Move the _GoBack bookmark from the end of the "Legal:" paragraph to the
#     (empty) ListParagraph that sits right before the "Data:" paragraph.
#  2. Rewrite the "Data:" bullet: "Data: 65 @ 15% = **9.75**" becomes
#     "Data: **21%**" (the paragraph mark/"next text" also becomes bold).
#  3. Rewrite the bold tail of the "HCI:" bullet: "**14%** **+** " becomes
#     "**33%** " (kept as two bold runs: "33%" and a trailing space).

$d = $word.ActiveDocument

function Get-ParagraphIndexByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Bookmark relocation: "_GoBack" moves from the "Legal:" paragraph to the
#    blank ListParagraph immediately before "Data:".
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$idxData = Get-ParagraphIndexByPrefix "Data:"
$idxBlankBeforeData = $idxData - 1
$pBlank = $d.Paragraphs.Item($idxBlankBeforeData)
$d.Bookmarks.Add("_GoBack", $pBlank.Range) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Data:" paragraph -> "Data: 21%" (bold "21%", bold paragraph mark).
# ---------------------------------------------------------------------------
$idxData = Get-ParagraphIndexByPrefix "Data:"
$pData = $d.Paragraphs.Item($idxData)

# Bold the whole paragraph first so the paragraph-mark run properties
# (<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>) get the bold flag too.
$pData.Range.Bold = 1

$pData = $d.Paragraphs.Item($idxData)
$dataStart = $pData.Range.Start
$dataEnd = $pData.Range.End
$rData = $d.Range($dataStart, $dataEnd - 1)

$dataXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
    '<w:r><w:t>Data:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>21%</w:t></w:r>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rData.InsertXML($dataXml)

# ---------------------------------------------------------------------------
# 3. "HCI:" paragraph -> bold tail "14% + " becomes "33% " (two bold runs).
# ---------------------------------------------------------------------------
$idxHCI = Get-ParagraphIndexByPrefix "HCI:"
$pHCI = $d.Paragraphs.Item($idxHCI)
$hciStart = $pHCI.Range.Start
$hciEnd = $pHCI.Range.End

# "HCI: " prefix is 5 characters; the rest (up to, excluding, the paragraph
# mark) is the bold "14% + " block that needs replacing.
$rHciBold = $d.Range($hciStart + 5, $hciEnd - 1)

$hciXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>33%</w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rHciBold.InsertXML($hciXml)
